# Update New Orleans xlsx: add a "State" column to hotel_info and
# reorder the worksheet tabs so review_info comes before hotel_info.

$wb = $excel.ActiveWorkbook
$wsHotel  = $wb.Worksheets.Item("hotel_info")
$wsReview = $wb.Worksheets.Item("review_info")

# Insert a new "State" column in hotel_info right after "Hotel_Name"
# (column A = STR, column B = Hotel_Name, so the new column is C,
# pushing City/Zip/... one to the right).
$wsHotel.Columns.Item(3).Insert()
$wsHotel.Range("C1").Value = "State"
$wsHotel.Range("C2").Value = "Louisiana"

# Reorder sheet tabs: move hotel_info to directly after review_info.
$wsHotel.Move($null, $wsReview)
